# Update the "Förändrad" (Changed) date column (C) from 2023-09-19 (45188)
# to 2023-09-20 (45189) for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
